# Auto-generated edit script: update FFXIV Leve profit-tracking values
# per external market-data scheduled-runner refresh (see commit message).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 129
$ws.Range("H129").Value = 907.2353000000001
$ws.Range("I129").Value = 261.75
$ws.Range("J129").Value = 993.3
$ws.Range("K129").Value = 785.25
$ws.Range("L129").Value = 2979.9
$ws.Range("M129").Value = 4214.75
$ws.Range("N129").Value = -12979.9
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 36000
$ws.Range("J136").Value = 36000
$ws.Range("L136").Value = 36000
$ws.Range("N136").Value = -46200
# Row 137
$ws.Range("H137").Value = 1793.8846
$ws.Range("I137").Value = 1575.762
$ws.Range("J137").Value = 2710
$ws.Range("K137").Value = 4727.286
$ws.Range("L137").Value = 8130
$ws.Range("M137").Value = -2177.286
$ws.Range("N137").Value = -13230
# Row 138
$ws.Range("H138").Value = 2152337.2
$ws.Range("I138").Value = 905.0714
$ws.Range("J138").Value = 3924105
$ws.Range("K138").Value = 2715.2142
$ws.Range("L138").Value = 11772315
$ws.Range("M138").Value = 2424.7858
$ws.Range("N138").Value = -11782595

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 977.4722
$ws.Range("I61").Value = 885
$ws.Range("J61").Value = 2549.5
$ws.Range("K61").Value = 885
$ws.Range("L61").Value = 2549.5
$ws.Range("M61").Value = -673
$ws.Range("N61").Value = -2973.5
# Row 62
$ws.Range("H62").Value = 29975
$ws.Range("J62").Value = 29975
$ws.Range("L62").Value = 29975
$ws.Range("N62").Value = -31223
# Row 65
$ws.Range("H65").Value = 29975
$ws.Range("J65").Value = 29975
$ws.Range("L65").Value = 89925
$ws.Range("N65").Value = -96165
# Row 74
$ws.Range("H74").Value = 30341.53
$ws.Range("I74").Value = 33042.324
$ws.Range("J74").Value = 2433.3333
$ws.Range("K74").Value = 33042.324
$ws.Range("L74").Value = 2433.3333
$ws.Range("M74").Value = -32168.324
$ws.Range("N74").Value = -4181.3333
# Row 77
$ws.Range("H77").Value = 30341.53
$ws.Range("I77").Value = 33042.324
$ws.Range("J77").Value = 2433.3333
$ws.Range("K77").Value = 165211.62
$ws.Range("L77").Value = 12166.6665
$ws.Range("M77").Value = -160843.62
$ws.Range("N77").Value = -20902.6665
# Row 132
$ws.Range("H132").Value = 1994.1897
$ws.Range("I132").Value = 1803.7838
$ws.Range("J132").Value = 2329.6667
$ws.Range("K132").Value = 5411.3514
$ws.Range("L132").Value = 6989.000100000001
$ws.Range("M132").Value = -2881.3514
$ws.Range("N132").Value = -12049.0001
# Row 136
$ws.Range("H136").Value = 977.4722
$ws.Range("I136").Value = 885
$ws.Range("J136").Value = 2549.5
$ws.Range("K136").Value = 2655
$ws.Range("L136").Value = 7648.5
$ws.Range("M136").Value = -105
$ws.Range("N136").Value = -12748.5
# Row 137
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200
# Row 138
$ws.Range("H138").Value = 38000
$ws.Range("J138").Value = 38000
$ws.Range("L138").Value = 38000
$ws.Range("N138").Value = -48280
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 6736.4546
$ws.Range("I25").Value = 1025.25
$ws.Range("K25").Value = 1025.25
$ws.Range("M25").Value = -790.25
# Row 57
$ws.Range("H57").Value = 45000
$ws.Range("J57").Value = 45000
$ws.Range("L57").Value = 45000
$ws.Range("N57").Value = -46440
# Row 81
$ws.Range("H81").Value = 34059.8
$ws.Range("J81").Value = 34059.8
$ws.Range("L81").Value = 34059.8
$ws.Range("N81").Value = -36181.8
# Row 84
$ws.Range("H84").Value = 34059.8
$ws.Range("J84").Value = 34059.8
$ws.Range("L84").Value = 102179.4
$ws.Range("N84").Value = -112787.4
# Row 134
$ws.Range("H134").Value = 4532.2446
$ws.Range("I134").Value = 3603.9143
$ws.Range("K134").Value = 10811.7429
$ws.Range("M134").Value = -8276.742899999999
# Row 136
$ws.Range("H136").Value = 45000
$ws.Range("J136").Value = 45000
$ws.Range("L136").Value = 45000
$ws.Range("N136").Value = -55200
# Row 137
$ws.Range("H137").Value = 39523.156
$ws.Range("J137").Value = 39523.156
$ws.Range("L137").Value = 39523.156
$ws.Range("N137").Value = -49723.156
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139").Value = 45700
$ws.Range("J139").Value = 43625
$ws.Range("L139").Value = 43625
$ws.Range("N139").Value = -53905
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
# Row 141
$ws.Range("H141").Value = 85333.336
$ws.Range("J141").Value = 85333.336
$ws.Range("L141").Value = 85333.336
$ws.Range("N141").Value = -95693.336

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 120
$ws.Range("I7").Value = 55
$ws.Range("J7").Value = 206.66667
$ws.Range("K7").Value = 55
$ws.Range("L7").Value = 206.66667
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = -432.66667
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 62
$ws.Range("H62").Value = 2852.6316
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 3612.5
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 3612.5
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -4860.5
# Row 65
$ws.Range("H65").Value = 2852.6316
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 3612.5
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 18062.5
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -24302.5
# Row 132
$ws.Range("H132").Value = 2053.2173
$ws.Range("I132").Value = 1371.1333
$ws.Range("J132").Value = 3332.125
$ws.Range("K132").Value = 4113.3999
$ws.Range("L132").Value = 9996.375
$ws.Range("M132").Value = -1583.3999
$ws.Range("N132").Value = -15056.375

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 30252
$ws.Range("J64").Value = 30252
$ws.Range("L64").Value = 30252
$ws.Range("N64").Value = -30702
# Row 67
$ws.Range("H67").Value = 30252
$ws.Range("J67").Value = 30252
$ws.Range("L67").Value = 30252
$ws.Range("N67").Value = -31812
# Row 68
$ws.Range("H68").Value = 20466.666
$ws.Range("I68").Value = 51500
$ws.Range("J68").Value = 4950
$ws.Range("K68").Value = 51500
$ws.Range("L68").Value = 4950
$ws.Range("M68").Value = -50751
$ws.Range("N68").Value = -6448
# Row 71
$ws.Range("H71").Value = 20466.666
$ws.Range("I71").Value = 51500
$ws.Range("J71").Value = 4950
$ws.Range("K71").Value = 257500
$ws.Range("L71").Value = 24750
$ws.Range("M71").Value = -253756
$ws.Range("N71").Value = -32238
# Row 132
$ws.Range("H132").Value = 2814.6316
$ws.Range("I132").Value = 2370.3076
$ws.Range("J132").Value = 3777.3333
$ws.Range("K132").Value = 7110.9228
$ws.Range("L132").Value = 11331.9999
$ws.Range("M132").Value = -4580.9228
$ws.Range("N132").Value = -16391.9999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 135
$ws.Range("H135").Value = 38000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140
